$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("K43").Value = 8748
$ws.Range("M43").Value = -8679
$ws.Range("I43").Value = 8748
$ws.Range("H43").Value = 13292.091
$ws.Range("K92").Value = 1128.6
$ws.Range("M92").Value = 119.4000000000001
$ws.Range("J92").Value = 1290.75
$ws.Range("L92").Value = 1290.75
$ws.Range("I92").Value = 1128.6
$ws.Range("H92").Value = 1155.625
$ws.Range("N92").Value = -3786.75
$ws.Range("K98").Value = 36079.812
$ws.Range("M98").Value = -34581.812
$ws.Range("I98").Value = 36079.812
$ws.Range("H98").Value = 32645.947
$ws.Range("K122").Value = 108239.436
$ws.Range("M122").Value = -105789.436
$ws.Range("I122").Value = 36079.812
$ws.Range("H122").Value = 32645.947
$ws.Range("K135").Value = 45867.1293
$ws.Range("M135").Value = -43332.1293
$ws.Range("I135").Value = 5096.3477
$ws.Range("H135").Value = 4688.357
$ws.Range("K138").Value = 1320611.82
$ws.Range("M138").Value = -1315471.82
$ws.Range("J138").Value = 5449.723
$ws.Range("L138").Value = 16349.169
$ws.Range("I138").Value = 440203.94
$ws.Range("H138").Value = 142740.53
$ws.Range("N138").Value = -26629.169

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("J5").Value = 2950
$ws.Range("L5").Value = 2950
$ws.Range("H5").Value = 4464.4287
$ws.Range("N5").Value = -3174
$ws.Range("K32").Value = 7208.5415
$ws.Range("M32").Value = -6921.5415
$ws.Range("J32").Value = 12500
$ws.Range("L32").Value = 12500
$ws.Range("I32").Value = 7208.5415
$ws.Range("H32").Value = 7487.0396
$ws.Range("N32").Value = -13074
$ws.Range("K45").Value = 4996.1763
$ws.Range("M45").Value = -4619.1763
$ws.Range("I45").Value = 4996.1763
$ws.Range("H45").Value = 5769.773
$ws.Range("K63").Value = 2166.1667
$ws.Range("M63").Value = -1480.1667
$ws.Range("I63").Value = 2166.1667
$ws.Range("H63").Value = 2571
$ws.Range("K66").Value = 10830.8335
$ws.Range("M66").Value = -7398.833500000001
$ws.Range("I66").Value = 2166.1667
$ws.Range("H66").Value = 2571
$ws.Range("K132").Value = 3827.5815
$ws.Range("M132").Value = -1297.5815
$ws.Range("I132").Value = 1275.8605
$ws.Range("H132").Value = 1955.0385

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("J4").Value = 2950
$ws.Range("L4").Value = 2950
$ws.Range("H4").Value = 4464.4287
$ws.Range("N4").Value = -3180
$ws.Range("K22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("H22").Value = 0
$ws.Range("N22").ClearContents()
$ws.Range("K86").Value = 7472.737
$ws.Range("M86").Value = -6349.737
$ws.Range("I86").Value = 7472.737
$ws.Range("H86").Value = 6844.5654
$ws.Range("K89").Value = 37363.685
$ws.Range("M89").Value = -31747.685
$ws.Range("I89").Value = 7472.737
$ws.Range("H89").Value = 6844.5654

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("K31").Value = 2018.5
$ws.Range("M31").Value = -1723.5
$ws.Range("I31").Value = 2018.5
$ws.Range("H31").Value = 3993.3572
$ws.Range("K34").Value = 2018.5
$ws.Range("M34").Value = -1816.5
$ws.Range("I34").Value = 2018.5
$ws.Range("H34").Value = 3993.3572
$ws.Range("K88").Value = 39155.5
$ws.Range("M88").Value = -38749.5
$ws.Range("J88").Value = 44662.668
$ws.Range("L88").Value = 44662.668
$ws.Range("I88").Value = 39155.5
$ws.Range("H88").Value = 42459.8
$ws.Range("N88").Value = -45474.668
$ws.Range("K91").Value = 39155.5
$ws.Range("M91").Value = -37751.5
$ws.Range("J91").Value = 44662.668
$ws.Range("L91").Value = 44662.668
$ws.Range("I91").Value = 39155.5
$ws.Range("H91").Value = 42459.8
$ws.Range("N91").Value = -47470.668
$ws.Range("K99").Value = 14760222
$ws.Range("M99").Value = -14758724
$ws.Range("I99").Value = 14760222
$ws.Range("H99").Value = 7874564
$ws.Range("K122").Value = 47818.875
$ws.Range("M122").Value = -45368.875
$ws.Range("I122").Value = 15939.625
$ws.Range("H122").Value = 12996.7
$ws.Range("K126").Value = 44280666
$ws.Range("M126").Value = -44278196
$ws.Range("I126").Value = 14760222
$ws.Range("H126").Value = 7874564
$ws.Range("J141").Value = 216406.39
$ws.Range("L141").Value = 216406.39
$ws.Range("H141").Value = 206655.83
$ws.Range("N141").Value = -226766.39

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("K17").Value = 1387.5
$ws.Range("M17").Value = -1218.5
$ws.Range("I17").Value = 462.5
$ws.Range("H17").Value = 1878.75
$ws.Range("J48").Value = 2013
$ws.Range("L48").Value = 6039
$ws.Range("H48").Value = 3103.9
$ws.Range("N48").Value = -6539
$ws.Range("K122").Value = 8772.300000000001
$ws.Range("M122").Value = -6322.300000000001
$ws.Range("I122").Value = 974.7
$ws.Range("H122").Value = 4492.636

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("K80").Value = 31079.8
$ws.Range("M80").Value = -30081.8
$ws.Range("I80").Value = 31079.8
$ws.Range("H80").Value = 18377.555
$ws.Range("K83").Value = 155399
$ws.Range("M83").Value = -150407
$ws.Range("I83").Value = 31079.8
$ws.Range("H83").Value = 18377.555
$ws.Range("J101").Value = 11000
$ws.Range("L101").Value = 11000
$ws.Range("H101").Value = 11000
$ws.Range("N101").Value = -17490
$ws.Range("K102").Value = 10198.357
$ws.Range("M102").Value = -8576.357
$ws.Range("I102").Value = 10198.357
$ws.Range("H102").Value = 9222.117
$ws.Range("J104").Value = 97666
$ws.Range("L104").Value = 97666
$ws.Range("H104").Value = 97666
$ws.Range("N104").Value = -104654
$ws.Range("J117").Value = 46000
$ws.Range("L117").Value = 46000
$ws.Range("H117").Value = 46000
$ws.Range("N117").Value = -52884
$ws.Range("K122").Value = 46631.142
$ws.Range("M122").Value = -44181.142
$ws.Range("J122").Value = 19874.25
$ws.Range("L122").Value = 59622.75
$ws.Range("I122").Value = 15543.714
$ws.Range("H122").Value = 17853.334
$ws.Range("N122").Value = -64522.75
$ws.Range("J123").Value = 40665
$ws.Range("L123").Value = 40665
$ws.Range("H123").Value = 40665
$ws.Range("N123").Value = -45565

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("K7").Value = 22459.268
$ws.Range("M7").Value = -22347.268
$ws.Range("I7").Value = 22459.268
$ws.Range("H7").Value = 19298.129
$ws.Range("K16").Value = 2196.818
$ws.Range("M16").Value = -2026.818
$ws.Range("J16").Value = 2756.1667
$ws.Range("L16").Value = 2756.1667
$ws.Range("I16").Value = 2196.818
$ws.Range("H16").Value = 2316.6785
$ws.Range("N16").Value = -3096.1667
$ws.Range("J40").Value = 17531.133
$ws.Range("L40").Value = 17531.133
$ws.Range("H40").Value = 29264
$ws.Range("N40").Value = -17803.133
$ws.Range("K46").Value = 787.375
$ws.Range("M46").Value = -599.375
$ws.Range("J46").Value = 2469.8572
$ws.Range("L46").Value = 2469.8572
$ws.Range("I46").Value = 787.375
$ws.Range("H46").Value = 1572.5333
$ws.Range("N46").Value = -2845.8572
$ws.Range("K100").Value = 2400
$ws.Range("M100").Value = -1859
$ws.Range("J100").Value = 7624.625
$ws.Range("L100").Value = 7624.625
$ws.Range("I100").Value = 2400
$ws.Range("H100").Value = 7044.1113
$ws.Range("N100").Value = -8706.625
$ws.Range("J104").Value = 30369
$ws.Range("L104").Value = 30369
$ws.Range("H104").Value = 30369
$ws.Range("N104").Value = -37357
$ws.Range("K118").Value = 25000
$ws.Range("M118").Value = -23343
$ws.Range("I118").Value = 25000
$ws.Range("H118").Value = 25000
$ws.Range("K122").Value = 18736.125
$ws.Range("M122").Value = -16286.125
$ws.Range("J122").Value = 6019.9287
$ws.Range("L122").Value = 18059.7861
$ws.Range("I122").Value = 6245.375
$ws.Range("H122").Value = 6140.1665
$ws.Range("N122").Value = -22959.7861
$ws.Range("K126").Value = 67377.804
$ws.Range("M126").Value = -64907.804
$ws.Range("I126").Value = 22459.268
$ws.Range("H126").Value = 19298.129
$ws.Range("K132").Value = 2811733.8
$ws.Range("M132").Value = -2809203.8
$ws.Range("J132").Value = 5575
$ws.Range("L132").Value = 16725
$ws.Range("I132").Value = 937244.6
$ws.Range("H132").Value = 626688.0600000001
$ws.Range("N132").Value = -21785

# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("J98").Value = 47999.75
$ws.Range("L98").Value = 47999.75
$ws.Range("H98").Value = 47999.75
$ws.Range("N98").Value = -53989.75
$ws.Range("J103").Value = 600
$ws.Range("L103").Value = 600
$ws.Range("H103").Value = 600
$ws.Range("N103").Value = -2944
$ws.Range("K122").Value = 6218.849999999999
$ws.Range("M122").Value = -3768.849999999999
$ws.Range("I122").Value = 2072.95
$ws.Range("H122").Value = 3970.0857
$ws.Range("K132").Value = 26350.431
$ws.Range("M132").Value = -23820.431
$ws.Range("J132").Value = 3349.6667
$ws.Range("L132").Value = 10049.0001
$ws.Range("I132").Value = 8783.477000000001
$ws.Range("H132").Value = 7353.5264
$ws.Range("N132").Value = -15109.0001
